$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A2 content from "coffee mug" to "Aleheida"
$ws.Range("A2").Value = "Aleheida"

# Update selection to A2
$ws.Range("A2").Select()
